$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each region block (Brasil: rows 2-22, Nordeste: rows 23-43, Sergipe: rows 44-64)
# has its Trimestre (C) / Valor (D) pairs shifted up by one row, with a new
# trailing quarter appended at the end of each block, and the Valor value cleared
# for rows that no longer have data.

# Keep the quarter labels stored as plain text (matching the original inline
# string cells) instead of being auto-converted to date serials.
$ws.Range("C2:C64").NumberFormat = "@"

$ws.Range("C2").Value = "01/10/2018"
$ws.Range("D2").Value = 88.28376452377132
$ws.Range("C3").Value = "01/01/2019"
$ws.Range("D3").Value = 87.15383963941923
$ws.Range("C4").Value = "01/04/2019"
$ws.Range("D4").Value = 87.85947559951479
$ws.Range("C5").Value = "01/07/2019"
$ws.Range("D5").Value = 88.09875854372994
$ws.Range("C6").Value = "01/10/2019"
$ws.Range("D6").Value = 88.91898936863468
$ws.Range("C7").Value = "01/01/2020"
$ws.Range("D7").Value = 87.6269256467444
$ws.Range("C8").Value = "01/04/2020"
$ws.Range("D8").Value = 86.40199837580568
$ws.Range("C9").Value = "01/07/2020"
$ws.Range("D9").Value = 85.10970347929863
$ws.Range("C10").Value = "01/10/2020"
$ws.Range("D10").Value = 85.82012456093744
$ws.Range("C11").Value = "01/01/2021"
$ws.Range("D11").Value = 85.09170501959174
$ws.Range("C12").Value = "01/04/2021"
$ws.Range("D12").Value = 85.7680202656022
$ws.Range("C13").Value = "01/07/2021"
$ws.Range("D13").Value = 87.35882739828995
$ws.Range("C14").Value = "01/10/2021"
$ws.Range("D14").Value = 88.85372779747212
$ws.Range("C15").Value = "01/01/2022"
$ws.Range("D15").Value = 88.85603969260613
$ws.Range("C16").Value = "01/04/2022"
$ws.Range("D16").Value = 90.6967300113522
$ws.Range("C17").Value = "01/07/2022"
$ws.Range("D17").Value = 91.29946932281176
$ws.Range("C18").Value = "01/10/2022"
$ws.Range("D18").Value = 92.05869818976858
$ws.Range("C19").Value = "01/01/2023"
$ws.Range("D19").Value = 91.20616836197172
$ws.Range("C20").Value = "01/04/2023"
$ws.Range("D20").Value = 91.96054185222719
$ws.Range("C21").Value = "01/07/2023"
$ws.Range("D21").Value = 92.31096399578379
$ws.Range("C22").Value = "01/10/2023"
$ws.Range("D22").Value = 92.59072488218143
$ws.Range("C23").Value = "01/10/2018"
$ws.Range("D23").Value = 85.46216809742464
$ws.Range("C24").Value = "01/01/2019"
$ws.Range("D24").Value = 84.55102285920174
$ws.Range("C25").Value = "01/04/2019"
$ws.Range("D25").Value = 85.24212449847059
$ws.Range("C26").Value = "01/07/2019"
$ws.Range("D26").Value = 85.39254559873116
$ws.Range("C27").Value = "01/10/2019"
$ws.Range("D27").Value = 86.24122476500217
$ws.Range("C28").Value = "01/01/2020"
$ws.Range("D28").Value = 84.21958837190678
$ws.Range("C29").Value = "01/04/2020"
$ws.Range("D29").ClearContents()
$ws.Range("C30").Value = "01/07/2020"
$ws.Range("D30").ClearContents()
$ws.Range("C31").Value = "01/10/2020"
$ws.Range("D31").ClearContents()
$ws.Range("C32").Value = "01/01/2021"
$ws.Range("D32").ClearContents()
$ws.Range("C33").Value = "01/04/2021"
$ws.Range("D33").ClearContents()
$ws.Range("C34").Value = "01/07/2021"
$ws.Range("D34").ClearContents()
$ws.Range("C35").Value = "01/10/2021"
$ws.Range("D35").ClearContents()
$ws.Range("C36").Value = "01/01/2022"
$ws.Range("D36").ClearContents()
$ws.Range("C37").Value = "01/04/2022"
$ws.Range("D37").Value = 87.30259251929546
$ws.Range("C38").Value = "01/07/2022"
$ws.Range("D38").Value = 88.02800283174703
$ws.Range("C39").Value = "01/10/2022"
$ws.Range("D39").Value = 89.13957176843775
$ws.Range("C40").Value = "01/01/2023"
$ws.Range("D40").Value = 87.76160329045526
$ws.Range("C41").Value = "01/04/2023"
$ws.Range("D41").Value = 88.67195362505535
$ws.Range("C42").Value = "01/07/2023"
$ws.Range("D42").Value = 89.1498039836851
$ws.Range("C43").Value = "01/10/2023"
$ws.Range("D43").Value = 89.56449309852451
$ws.Range("C44").Value = "01/10/2018"
$ws.Range("D44").Value = 84.78873239436619
$ws.Range("C45").Value = "01/01/2019"
$ws.Range("D45").Value = 84.47789275634995
$ws.Range("C46").Value = "01/04/2019"
$ws.Range("D46").Value = 84.6503178928247
$ws.Range("C47").Value = "01/07/2019"
$ws.Range("D47").Value = 85.21897810218978
$ws.Range("C48").Value = "01/10/2019"
$ws.Range("D48").Value = 85.06666666666666
$ws.Range("C49").Value = "01/01/2020"
$ws.Range("D49").Value = 84.21052631578947
$ws.Range("C50").Value = "01/04/2020"
$ws.Range("D50").ClearContents()
$ws.Range("C51").Value = "01/07/2020"
$ws.Range("D51").ClearContents()
$ws.Range("C52").Value = "01/10/2020"
$ws.Range("D52").ClearContents()
$ws.Range("C53").Value = "01/01/2021"
$ws.Range("D53").ClearContents()
$ws.Range("C54").Value = "01/04/2021"
$ws.Range("D54").ClearContents()
$ws.Range("C55").Value = "01/07/2021"
$ws.Range("D55").ClearContents()
$ws.Range("C56").Value = "01/10/2021"
$ws.Range("D56").ClearContents()
$ws.Range("C57").Value = "01/01/2022"
$ws.Range("D57").ClearContents()
$ws.Range("C58").Value = "01/04/2022"
$ws.Range("D58").Value = 87.24954462659382
$ws.Range("C59").Value = "01/07/2022"
$ws.Range("D59").Value = 87.87037037037037
$ws.Range("C60").Value = "01/10/2022"
$ws.Range("D60").Value = 88.04744525547446
$ws.Range("C61").Value = "01/01/2023"
$ws.Range("D61").Value = 88.1740775780511
$ws.Range("C62").Value = "01/04/2023"
$ws.Range("D62").Value = 89.76303317535546
$ws.Range("C63").Value = "01/07/2023"
$ws.Range("D63").Value = 90.20332717190388
$ws.Range("C64").Value = "01/10/2023"
$ws.Range("D64").Value = 88.70214752567693
